# use K+ for combined dropped third strike
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (Droppped Third Strike): add VAR3 = "K+"
$ws.Range("F31").Value = "K+"

# Row 32 (Passed Ball on Dropped Third Strike): VAR1 becomes "Kd" -> "PB", VAR3 becomes "PB" -> "K+"
$ws.Range("D32").Value = "PB"
$ws.Range("F32").Value = "K+"

# Row 33 (Wild Pitch on Dropped Third Strike): VAR1 becomes "Kd" -> "WP", VAR3 becomes "WP" -> "K+"
$ws.Range("D33").Value = "WP"
$ws.Range("F33").Value = "K+"

# Row 43 (Passed Ball, Advance): VAR1 becomes "E" -> "PB", VAR3 cleared
$ws.Range("D43").Value = "PB"
$ws.Range("F43").Clear()

# Row 46 (Wild Pitch, Advance): VAR1 becomes "E" -> "WP", VAR3 cleared
$ws.Range("D46").Value = "WP"
$ws.Range("F46").Clear()

# Update view state to match authored selection/scroll position
$ws.Range("E28").Select() | Out-Null
